$d = $word.ActiveDocument

# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph and the
# "© 2020 . Contact: ..." paragraph entirely (including their paragraph marks),
# while leaving the surrounding empty paragraphs untouched.

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $p.Range.Delete()
    }
    elseif ($t -like "*Contact: luizeleno@usp.br*") {
        $p.Range.Delete()
    }
}
